# Weekly update: a new daily price record for "Piña" (pineapple) is added
# for the Feria Lagunitas de Puerto Montt market. Insert a new row at 53
# (pushing the existing row 53 and everything below it down by one) and
# populate it with the new day's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(53).Insert()

$ws.Cells.Item(53,1).Value = 4
$ws.Cells.Item(53,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(53,3).Value = "Los Lagos"
$ws.Cells.Item(53,4).Value = 44544
$ws.Cells.Item(53,5).Value = 10
$ws.Cells.Item(53,6).Value = "Fruta"
$ws.Cells.Item(53,7).Value = 100108
$ws.Cells.Item(53,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(53,9).Value = 100108005
$ws.Cells.Item(53,10).Value = "Piña"
$ws.Cells.Item(53,11).Value = "Caramelo"
$ws.Cells.Item(53,12).Value = "Tercera"
$ws.Cells.Item(53,13).Value = 200
$ws.Cells.Item(53,14).Value = 20000
$ws.Cells.Item(53,15).Value = 21000
$ws.Cells.Item(53,16).Value = 20500
$ws.Cells.Item(53,17).Value = '$/caja 16 unidades'
$ws.Cells.Item(53,18).Value = "Ecuador"
$ws.Cells.Item(53,19).Value = 1281
$ws.Cells.Item(53,20).Value = 16
